$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 41/42: Coin name and Link swap (plain text columns, safe to assign directly) ---
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"

# --- Price (D) and Volume(1h) (E) columns store numeric-looking strings as text. ---
# --- Force text format before assigning so Excel does not coerce them to numbers/percentages, ---
# --- then restore the original "General" number format afterwards. ---
$priceVolumeRange = $ws.Range("D2:E50")
$priceVolumeRange.NumberFormat = "@"

$ws.Range("D2").Value = "257.71"
$ws.Range("E2").Value = "0.59%"
$ws.Range("D3").Value = "27.13"
$ws.Range("E3").Value = "-4.17%"
$ws.Range("D4").Value = "4.915"
$ws.Range("E4").Value = "-7.41%"
$ws.Range("D5").Value = "0.05950"
$ws.Range("E5").Value = "2.62%"
$ws.Range("D6").Value = "6.693"
$ws.Range("E6").Value = "-0.04%"
$ws.Range("D7").Value = "0.8688"
$ws.Range("E7").Value = "-0.30%"
$ws.Range("D8").Value = "0.9559"
$ws.Range("E8").Value = "4.68%"
$ws.Range("D9").Value = "0.1411"
$ws.Range("E9").Value = "0.16%"
$ws.Range("D10").Value = "0.03542"
$ws.Range("E10").Value = "3.89%"
$ws.Range("D11").Value = "0.07174"
$ws.Range("E11").Value = "0.13%"
$ws.Range("D12").Value = "0.03144"
$ws.Range("E12").Value = "-0.48%"
$ws.Range("D13").Value = "0.09255"
$ws.Range("E13").Value = "0.33%"
$ws.Range("D14").Value = "0.001551"
$ws.Range("E14").Value = "-0.61%"
$ws.Range("D15").Value = "0.0006060"
$ws.Range("E15").Value = "0.36%"
$ws.Range("D16").Value = "0.006013"
$ws.Range("E16").Value = "1.75%"
$ws.Range("D17").Value = "3.485"
$ws.Range("E17").Value = "-0.69%"
$ws.Range("D18").Value = "3.258"
$ws.Range("E18").Value = "0.72%"
$ws.Range("E19").Value = "-2.87%"
$ws.Range("E20").Value = "0.59%"
$ws.Range("E21").Value = "-0.46%"
$ws.Range("D22").Value = "3.533"
$ws.Range("E22").Value = "0.43%"
$ws.Range("D23").Value = "0.04284"
$ws.Range("E23").Value = "2.72%"
$ws.Range("E24").Value = "2.59%"
$ws.Range("D25").Value = "0.001220"
$ws.Range("E25").Value = "-0.12%"
$ws.Range("D26").Value = "0.004519"
$ws.Range("E26").Value = "-9.22%"
$ws.Range("E27").Value = "0.20%"
$ws.Range("E28").Value = "-22.92%"
$ws.Range("D40").Value = "0.03833"
$ws.Range("E40").Value = "-0.96%"
$ws.Range("D41").Value = "0.1103"
$ws.Range("E41").Value = "0.57%"
$ws.Range("D42").Value = "0.003984"
$ws.Range("E42").Value = "-30.29%"
$ws.Range("D43").Value = "0.002310"
$ws.Range("E43").Value = "-0.05%"
$ws.Range("E44").Value = "0.11%"
$ws.Range("D45").Value = "0.00005490"
$ws.Range("E45").Value = "4.21%"
$ws.Range("E46").Value = "0.19%"
$ws.Range("E47").Value = "28.59%"
$ws.Range("D48").Value = "0.002130"
$ws.Range("E48").Value = "-1.81%"
$ws.Range("E49").Value = "0.19%"
$ws.Range("E50").Value = "0.19%"

$priceVolumeRange.NumberFormat = "General"
